# Apply the data edits described by the commit "modified data and added graphs/correlations"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Several existing "told" responses in column D were re-coded to "seen"
# ---------------------------------------------------------------------------
$ws.Range("D2").Value  = "seen"
$ws.Range("D4").Value  = "seen"
$ws.Range("D12").Value = "seen"
$ws.Range("D17").Value = "seen"
$ws.Range("D19").Value = "seen"

# Row 21 used the capitalized "Told" - normalize to lowercase "told"
$ws.Range("D21").Value = "told"

# ---------------------------------------------------------------------------
# 2. Fill in the previously-blank rows 47-51 with new survey data, and add a
#    brand new row 49. Rows 47, 48, 50 and 51 already carried their shading /
#    number-format styling on column A - only the values need to be written.
# ---------------------------------------------------------------------------
$ws.Range("A47").Value = "Male"
$ws.Range("B47").Value = 48
$ws.Range("C47").Value = 66
$ws.Range("D47").Value = "told"

$ws.Range("A48").Value = "Female"
$ws.Range("B48").Value = 45
$ws.Range("C48").Value = 29
$ws.Range("D48").Value = "told"

# Row 49 did not previously exist - give A49 the same "0.00" number format
# used by the other rows in this block (style id 2) before filling it in.
$ws.Range("A49").NumberFormat = "0.00"
$ws.Range("A49").Value = "Female"
$ws.Range("B49").Value = 20
$ws.Range("C49").Value = 83
$ws.Range("D49").Value = "told"

$ws.Range("A50").Value = "Male"
$ws.Range("B50").Value = 19
$ws.Range("C50").Value = 132
$ws.Range("D50").Value = "seen"

$ws.Range("A51").Value = "Male"
$ws.Range("B51").Value = 21
$ws.Range("C51").Value = 89
$ws.Range("D51").Value = "seen"

# ---------------------------------------------------------------------------
# 3. Update the saved view: scroll back up to the top of the sheet and
#    reselect D2 (previously the window was scrolled to row 28 with D40
#    selected).
# ---------------------------------------------------------------------------
$ws.Range("D2").Select()
